$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 13008.556
$ws.Range("I2").Value = 2512.5
$ws.Range("J2").Value = 34000.668
$ws.Range("K2").Value = 2512.5
$ws.Range("L2").Value = 34000.668
$ws.Range("M2").Value = -2399.5
$ws.Range("N2").Value = -34226.668
$ws.Range("H9").Value = 1000155.7
$ws.Range("I9").Value = 279.3
$ws.Range("J9").Value = 1714353.1
$ws.Range("K9").Value = 279.3
$ws.Range("L9").Value = 1714353.1
$ws.Range("M9").Value = -110.3
$ws.Range("N9").Value = -1714691.1
$ws.Range("H40").Value = 3553.5557
$ws.Range("I40").Value = 3098.04
$ws.Range("K40").Value = 3098.04
$ws.Range("M40").Value = -2923.04
$ws.Range("H58").Value = 28211.334
$ws.Range("I58").Value = 7700.2856
$ws.Range("K58").Value = 23100.8568
$ws.Range("M58").Value = -22950.8568
$ws.Range("H61").Value = 141.66667
$ws.Range("I61").Value = 83.333336
$ws.Range("K61").Value = 250.000008
$ws.Range("M61").Value = -78.00000800000001
$ws.Range("H64").Value = 4620.0625
$ws.Range("I64").Value = 4250.125
$ws.Range("J64").Value = 4990
$ws.Range("K64").Value = 4250.125
$ws.Range("L64").Value = 4990
$ws.Range("M64").Value = -4002.125
$ws.Range("N64").Value = -5486
$ws.Range("H67").Value = 4620.0625
$ws.Range("I67").Value = 4250.125
$ws.Range("J67").Value = 4990
$ws.Range("K67").Value = 4250.125
$ws.Range("L67").Value = 4990
$ws.Range("M67").Value = -3392.125
$ws.Range("N67").Value = -6706
$ws.Range("H93").Value = 25000
$ws.Range("J93").Value = 25000
$ws.Range("L93").Value = 25000
$ws.Range("N93").Value = -29992
$ws.Range("H113").Value = 252204.38
$ws.Range("I113").Value = 2937.1667
$ws.Range("J113").Value = 1000006
$ws.Range("K113").Value = 2937.1667
$ws.Range("L113").Value = 1000006
$ws.Range("M113").Value = 316.8332999999998
$ws.Range("N113").Value = -1006514
$ws.Range("H132").Value = 1572.2354
$ws.Range("I132").Value = 1382
$ws.Range("K132").Value = 4146
$ws.Range("M132").Value = -1616
$ws.Range("H138").Value = 3781.549
$ws.Range("I138").Value = 1476.6522
$ws.Range("J138").Value = 5674.857
$ws.Range("K138").Value = 4429.9566
$ws.Range("L138").Value = 17024.571
$ws.Range("M138").Value = 710.0434000000005
$ws.Range("N138").Value = -27304.571

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 6879.905
$ws.Range("I45").Value = 7562.294
$ws.Range("J45").Value = 3979.75
$ws.Range("K45").Value = 7562.294
$ws.Range("L45").Value = 3979.75
$ws.Range("M45").Value = -7185.294
$ws.Range("N45").Value = -4733.75
$ws.Range("H74").Value = 567078.5600000001
$ws.Range("I74").Value = 1007791.5
$ws.Range("K74").Value = 1007791.5
$ws.Range("M74").Value = -1006917.5
$ws.Range("H77").Value = 567078.5600000001
$ws.Range("I77").Value = 1007791.5
$ws.Range("K77").Value = 5038957.5
$ws.Range("M77").Value = -5034589.5
$ws.Range("H98").Value = 57333
$ws.Range("J98").Value = 57333
$ws.Range("L98").Value = 57333
$ws.Range("N98").Value = -63323

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H33").Value = 0
$ws.Range("I33").Value = 0
$ws.Range("K33").Value = 0
$ws.Range("M33").ClearContents()
$ws.Range("H94").Value = 930.9048
$ws.Range("I94").Value = 1003.25
$ws.Range("J94").Value = 699.4
$ws.Range("K94").Value = 1003.25
$ws.Range("L94").Value = 699.4
$ws.Range("M94").Value = -552.25
$ws.Range("N94").Value = -1601.4
$ws.Range("H105").Value = 50014508
$ws.Range("I105").Value = 66684980
$ws.Range("K105").Value = 66684980
$ws.Range("M105").Value = -66683233
$ws.Range("H134").Value = 3791.611
$ws.Range("I134").Value = 1016.3
$ws.Range("K134").Value = 3048.9
$ws.Range("M134").Value = -513.8999999999996

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H28").Value = 34234.57
$ws.Range("J28").Value = 34234.57
$ws.Range("L28").Value = 34234.57
$ws.Range("N28").Value = -34724.57
$ws.Range("H56").Value = 20925
$ws.Range("J56").Value = 20925
$ws.Range("L56").Value = 20925
$ws.Range("N56").Value = -22615
$ws.Range("H60").Value = 26499.727
$ws.Range("I60").Value = 7916.6665
$ws.Range("J60").Value = 48799.4
$ws.Range("K60").Value = 7916.6665
$ws.Range("L60").Value = 48799.4
$ws.Range("M60").Value = -7405.6665
$ws.Range("N60").Value = -49821.4
$ws.Range("H63").Value = 0
$ws.Range("J63").Value = 0
$ws.Range("L63").Value = 0
$ws.Range("N63").ClearContents()
$ws.Range("H66").Value = 0
$ws.Range("J66").Value = 0
$ws.Range("L66").Value = 0
$ws.Range("N66").ClearContents()

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 1546.2222
$ws.Range("I113").Value = 1196.5454
$ws.Range("K113").Value = 3589.6362
$ws.Range("M113").Value = -1419.6362
$ws.Range("H131").Value = 3717.3438
$ws.Range("I131").Value = 1074.6666
$ws.Range("J131").Value = 4751.4346
$ws.Range("K131").Value = 3223.9998
$ws.Range("L131").Value = 14254.3038
$ws.Range("M131").Value = 1816.0002
$ws.Range("N131").Value = -24334.3038

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H17").Value = 5497.222
$ws.Range("J17").Value = 6096.875
$ws.Range("L17").Value = 6096.875
$ws.Range("N17").Value = -6432.875
$ws.Range("H32").Value = 0
$ws.Range("I32").Value = 0
$ws.Range("J32").Value = 0
$ws.Range("K32").Value = 0
$ws.Range("L32").Value = 0
$ws.Range("M32").ClearContents()
$ws.Range("N32").ClearContents()
$ws.Range("H45").Value = 0
$ws.Range("J45").Value = 0
$ws.Range("L45").Value = 0
$ws.Range("N45").ClearContents()
$ws.Range("H70").Value = 14213
$ws.Range("J70").Value = 16698.2
$ws.Range("L70").Value = 16698.2
$ws.Range("N70").Value = -17238.2
$ws.Range("H73").Value = 14213
$ws.Range("J73").Value = 16698.2
$ws.Range("L73").Value = 16698.2
$ws.Range("N73").Value = -18570.2
$ws.Range("H102").Value = 4484.4546
$ws.Range("I102").Value = 4631.5
$ws.Range("K102").Value = 4631.5
$ws.Range("M102").Value = -3009.5
$ws.Range("H107").Value = 629.375
$ws.Range("I107").Value = 506
$ws.Range("J107").Value = 999.5
$ws.Range("K107").Value = 506
$ws.Range("L107").Value = 999.5
$ws.Range("M107").Value = 1414
$ws.Range("N107").Value = -4839.5
$ws.Range("H126").Value = 3169.625
$ws.Range("I126").Value = 2407.8
$ws.Range("K126").Value = 7223.400000000001
$ws.Range("M126").Value = -4753.400000000001
$ws.Range("H132").Value = 5702.913
$ws.Range("I132").Value = 4408.6
$ws.Range("K132").Value = 13225.8
$ws.Range("M132").Value = -10695.8

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 80568.766
$ws.Range("I7").Value = 94349.27
$ws.Range("K7").Value = 94349.27
$ws.Range("M7").Value = -94237.27
$ws.Range("H9").Value = 733.7778
$ws.Range("I9").Value = 202.66667
$ws.Range("J9").Value = 999.3333
$ws.Range("K9").Value = 202.66667
$ws.Range("L9").Value = 999.3333
$ws.Range("M9").Value = 21.33332999999999
$ws.Range("N9").Value = -1447.3333
$ws.Range("H16").Value = 802.44446
$ws.Range("I16").Value = 873.4783
$ws.Range("K16").Value = 873.4783
$ws.Range("M16").Value = -703.4783
$ws.Range("H46").Value = 6175.4194
$ws.Range("I46").Value = 1401
$ws.Range("J46").Value = 6334.567
$ws.Range("K46").Value = 1401
$ws.Range("L46").Value = 6334.567
$ws.Range("M46").Value = -1213
$ws.Range("N46").Value = -6710.567
$ws.Range("H126").Value = 80568.766
$ws.Range("I126").Value = 94349.27
$ws.Range("K126").Value = 283047.81
$ws.Range("M126").Value = -280577.81
$ws.Range("H132").Value = 4810.485
$ws.Range("I132").Value = 4049.2144
$ws.Range("K132").Value = 12147.6432
$ws.Range("M132").Value = -9617.643199999999

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H2").Value = 200617.8
$ws.Range("I2").Value = 44.5
$ws.Range("J2").Value = 334333.34
$ws.Range("K2").Value = 44.5
$ws.Range("L2").Value = 334333.34
$ws.Range("M2").Value = 67.5
$ws.Range("N2").Value = -334557.34
$ws.Range("H8").Value = 125680
$ws.Range("I8").Value = 143562.86
$ws.Range("K8").Value = 143562.86
$ws.Range("M8").Value = -143422.86
$ws.Range("H63").Value = 44500
$ws.Range("J63").Value = 44500
$ws.Range("L63").Value = 44500
$ws.Range("N63").Value = -45748
$ws.Range("H66").Value = 44500
$ws.Range("J66").Value = 44500
$ws.Range("L66").Value = 133500
$ws.Range("N66").Value = -139740
$ws.Range("H123").Value = 0
$ws.Range("J123").Value = 0
$ws.Range("L123").Value = 0
$ws.Range("N123").ClearContents()
$ws.Range("H129").Value = 47250
$ws.Range("J129").Value = 47250
$ws.Range("L129").Value = 47250
$ws.Range("N129").Value = -57250
